# Apply the "Gastos" workbook edit described in the commit:
# "Adicionei uma pasta com o projeto do pupeetter"
#
# Functional changes to the worksheet data:
#  - C24 updated from 504 to 512
#  - C26, C27, C28, C29, C30 cleared (contents removed)
#  - C40 updated from 1735.22 to 541.44000000000005
#  - C41, C42 cleared (contents removed)
#  - C51, D51, E51, G51, H51 cleared (contents removed)
#  - Active cell/selection left on F3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 - value bump
$ws.Range("C24").Value = 512

# Rows 26-30 - clear out the first-account column values
$ws.Range("C26:C30").ClearContents()

# Row 40 - update hard-coded starting balance
$ws.Range("C40").Value = 541.44000000000005

# Rows 41-42 - clear extra entries
$ws.Range("C41:C42").ClearContents()

# Row 51 - clear the whole entry (values + the SUM formula in G51)
$ws.Range("C51:E51").ClearContents()
$ws.Range("G51:H51").ClearContents()

# Leave the cursor on F3, matching the saved selection in the workbook
$ws.Range("F3").Select()
